# Release: bump splash screen date placeholder (slide master + all slide
# layouts) from 8/3/2010 -> 8/17/2010, and the "Release 0.3.0" text on the
# splash slide to "Release 1.0.0".

$p = $ppt.ActivePresentation

# --- Update the "datetimeFigureOut" Date Placeholder text wherever it
# appears: the slide master and every custom (slide) layout. ---
$sm = $p.SlideMaster

for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $shp = $sm.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq "8/3/2010") {
            $shp.TextFrame.TextRange.Text = "8/17/2010"
        }
    }
}

for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    for ($si = 1; $si -le $cl.Shapes.Count; $si++) {
        $shp = $cl.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "8/3/2010") {
                $shp.TextFrame.TextRange.Text = "8/17/2010"
            }
        }
    }
}

# --- Update the splash-slide "Release 0.3.0" text to "Release 1.0.0",
# only touching the version-number run so the "Release " run is left
# untouched. ---
$s = $p.Slides.Item(1)
for ($si = 1; $si -le $s.Shapes.Count; $si++) {
    $shp = $s.Shapes.Item($si)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Release 0.3.0") {
            $sub = $tr.Characters(9, 5)
            $sub.Text = "1.0.0"
        }
    }
}
